$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the text/filename/function cells of row 4 entirely (value + formatting)
$ws.Range("A4:C4").Clear()
# Clear only the numeric value in D4, keeping its cell style
$ws.Range("D4").ClearContents()

# Move the active selection to C5, matching the saved selection in the sheet
$ws.Range("C5").Select()
